$d = $word.ActiveDocument

function Get-ParagraphIndexByText {
    param($doc, $needle)
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $needle) {
            return $i
        }
    }
    return -1
}

function Add-ItalicParagraphAfter {
    param($doc, $afterIndex, $text)
    $afterPar = $doc.Paragraphs.Item($afterIndex)
    $null = $afterPar.Range.InsertParagraphAfter()
    $newPar = $doc.Paragraphs.Item($afterIndex + 1)
    $start = $newPar.Range.Start
    $insPoint = $doc.Range($start, $start)
    $insPoint.InsertAfter($text)
    $contentRange = $doc.Range($start, $start + $text.Length)
    $contentRange.Font.Italic = $true
}

# -----------------------------------------------------------------
# Text blocks (verbatim from the diff)
# -----------------------------------------------------------------
$ativacao_old = "Ativação: 01/01/2012"
$ativacao_new = "Ativação: 01/01/2022"

$objetivos_pt_old = "Fornecer uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base, de transformação."
$objetivos_pt_new = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."
$objetivos_en_new = "Provide students with a current view of industrial processes that use chemical conversion as a route to transform raw material into product. The processes of the basic chemical and transformation industries will be studied."

$resumido_pt_old = "Introdução ao Estudo dos Processos Químicos Industriais. Relacionamento com a Engenharia Química; Derivados Inorgânicos do Nitrogênio; Ácido Sulfúrico; Fósforo e Ácido Fosfórico; Fertilizantes; Indústrias de Cloro Álcalis; Indústrias de Vidros e Materiais Cerâmicos; Carga e Pigmentos Inorgânicos."
$resumido_pt_new = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"
$resumido_en_new = "Introduction to Industrial Chemical Processes; NPK / Fertilizers; Sulfuric Acid; Chlorine Alkali; Paper and Cellulose; Sugar and alcohol; Biotechnological Processes."

$programa_pt_old = "Introdução ao estudo dos Processos Químicos Industriais. Relacionamento com a Engenharia Química. Fundamentos dos processos químicos. Condução dos processos (batelada X contínuo). Fluxogramas. Derivados inorgânicos do nitrogênio - Introdução Amônia. Generalidades. Amônia. Produção sintética pelo processo Haber Bosch. Uréia: Generalidades. Processo de Fabricação. Nitrato de Amônia: Generalidades - Processo de Fabricação. Acído Nítrico. Generalidades. Processo de Fabricação Ácido Sulfúrico. Generalidades. Processo de Fabricação. Concentração. Fósforo e Ácido fosfórico. Generalidades. Matérias Primas. Produção de ácido fosfórico.  Indústrias de cloro álcalis. Generalidades. Matérias primas. Produção de barrilha e bicarbonato de sódio. Indústria de cloro e álcalis: produção de cloro e soda caústica. Células a diafragma.e mercúrio. Comparação. Ácido Clorídrico: fabricação e aplicações. Indústrias de vidros e materiais cerâmicos: matérias primas e fabricação."
$programa_pt_new = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."
$programa_en_new = "1. Introduction to Industrial Chemical Processes;2. NPK / Fertilizers3. Sulfuric Acid;4. Chlorine Alkali;5. Paper and Cellulose;6. Sugar and alcohol;7. Biotechnological Processes;"

$metodo_old = "Serão aplicadas duas provas (P1 e P2) e a Nota Final (NF) será a média aritmética das mesmas."
$metodo_new = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."

$criterio_old = "Serão aprovados os alunos com NF maior ou igual a 5,0 e frequência superior a 70%."
$criterio_new = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."

$norma_old = "Será feita a Recuperação( REC) para alunos com NF maior ou igual a 3,0 e menor que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou superior a 5,0, sendo MF= (NF+ REC)/2 .                                                 Na semana anterior à REC será dada uma aula de recordação de toda a matéria apresentada."
$norma_new = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

$bib_new = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"

# -----------------------------------------------------------------
# 1) Ativação: 01/01/2012 -> 01/01/2022
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute($ativacao_old, $false, $false, $false, $false, $false, $true, 1, $false, $ativacao_new, 2)

# -----------------------------------------------------------------
# 2) Objetivos: replace PT text, add italic EN paragraph after
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute($objetivos_pt_old, $false, $false, $false, $false, $false, $true, 1, $false, $objetivos_pt_new, 2)
$idx = Get-ParagraphIndexByText $d $objetivos_pt_new
Add-ItalicParagraphAfter $d $idx $objetivos_en_new

# -----------------------------------------------------------------
# 3) Programa resumido: replace PT text, add italic EN paragraph after
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute($resumido_pt_old, $false, $false, $false, $false, $false, $true, 1, $false, $resumido_pt_new, 2)
$idx = Get-ParagraphIndexByText $d $resumido_pt_new
Add-ItalicParagraphAfter $d $idx $resumido_en_new

# -----------------------------------------------------------------
# 4) Programa: replace PT text, add italic EN paragraph after
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute($programa_pt_old, $false, $false, $false, $false, $false, $true, 1, $false, $programa_pt_new, 2)
$idx = Get-ParagraphIndexByText $d $programa_pt_new
Add-ItalicParagraphAfter $d $idx $programa_en_new

# -----------------------------------------------------------------
# 5) Avaliação: Método / Critério / Norma de recuperação values
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute($metodo_old, $false, $false, $false, $false, $false, $true, 1, $false, $metodo_new, 2)
$null = $d.Content.Find.Execute($criterio_old, $false, $false, $false, $false, $false, $true, 1, $false, $criterio_new, 2)
$null = $d.Content.Find.Execute($norma_old, $false, $false, $false, $false, $false, $true, 1, $false, $norma_new, 2)

# -----------------------------------------------------------------
# 6) Bibliografia: replace whole paragraph (9 runs + breaks) with a
#    single clean run (no breaks, no xml:space="preserve").
# -----------------------------------------------------------------
$bibIdx = Get-ParagraphIndexByText $d "Bibliografia"
$bibBodyPar = $d.Paragraphs.Item($bibIdx + 1)
$bibRange = $bibBodyPar.Range
$bibRange.End = $bibRange.End - 1
$bibRange.Text = ""
$bibBodyPar2 = $d.Paragraphs.Item($bibIdx + 1)
$bibRange2 = $bibBodyPar2.Range
$bibRange2.End = $bibRange2.End - 1
$bibRange2.InsertAfter($bib_new)

Write-Host "All edits applied"
